# Update Excel files after daily scrape - 2025-07-27 03:51:13 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes ---
# Excel's ColumnWidth setter pads stored width by 5/6 (0.8333...) relative to
# the value that ends up written to the OOXML <col width="..."/>. Subtract
# that offset so the persisted width matches the target exactly.
$padding = 5 / 6
$ws.Columns.Item(3).ColumnWidth = 37 - $padding
$ws.Columns.Item(4).ColumnWidth = 22 - $padding
$ws.Columns.Item(6).ColumnWidth = 15 - $padding
$ws.Columns.Item(8).ColumnWidth = 27 - $padding

# --- Row 2 updates ---
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1326581"
$ws.Range("C2").Value = "Culture & Development Responsible"
$ws.Range("D2").Value = "Zurique, Suíça"
$ws.Range("F2").Value = "7 applicants"
$ws.Range("H2").Value = "SIX Group Services AG"

# A2 holds a numeric-looking id that must remain stored as text, exactly
# like the original inline string. Force the cell to Text format first so
# the engine keeps it as a string instead of re-typing it as a number.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1326581"

# --- Row 3 updates ---
$ws.Range("C3").Value = "Sales & Event Executive (3 months)"
$ws.Range("D3").Value = "Nugegoda, Sri Lanka"
$ws.Range("F3").Value = "1 applicant"
$ws.Range("G3").Value = "3 - 6 Months"
$ws.Range("H3").Value = "Brand Corridor (Pvt) Ltd"

$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1325142"
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "1325142"

# --- Remove rows 4-7 (only 2 data rows remain after this edit) ---
$ws.Range("A4:H7").EntireRow.Delete()
